$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New sprint-log entries added at the bottom of the 4th sprint table.
# The shared-string table is append-on-first-use, so the cells must be
# written in the exact order that reproduces the target shared string
# indices (40..45):
#   C41 -> "Vorbereitung für Präsentation"                                   (idx 40)
#   C44 -> "Layout für Tablets erstellt"                                     (idx 41)
#   C42 -> "Kleinere Bugfixes"                                               (idx 42)
#   C43 -> "Besprechung bzgl. Präsentration"                                 (idx 43)
#   C45 -> "Letzte Bugfixes für Präsentation"                                (idx 44)
#   C46 -> "Besprechung des neuen Sprints mithilfe der neu gewonnen Ideen der Präsentationen" (idx 45)
$ws.Range("C41").Value = "Vorbereitung für Präsentation"
$ws.Range("C44").Value = "Layout für Tablets erstellt"
$ws.Range("C42").Value = "Kleinere Bugfixes"
$ws.Range("C43").Value = "Besprechung bzgl. Präsentration"
$ws.Range("C45").Value = "Letzte Bugfixes für Präsentation"
$ws.Range("C46").Value = "Besprechung des neuen Sprints mithilfe der neu gewonnen Ideen der Präsentationen"

# C47 becomes a "Frei!" marker, matching the styling/text already used for
# the same marker elsewhere in the sheet (e.g. C27).
$ws.Range("C27").Copy()
$ws.Range("C47").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C47").Value = $ws.Range("C27").Value2

# A new, empty, formatted cell is appended further down the sheet.
$ws.Range("C27").Copy()
$ws.Range("F51").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the visible selection to match the saved view state.
[void]$ws.Range("F49").Select()
